$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Edit 1: "Data di consegna: 24 ottobre 2025" -> "Data di consegna: 31 ottobre 2025"
# The final document splits this text into three runs:
#   "Data di consegna: " | "31" | " ottobre 2025"
# Locate the exact character offset of "24" dynamically (no hard-coded
# positions) by searching for the stable prefix text first.
# -------------------------------------------------------------------------
$prefixRng = $d.Content
$prefixFound = $prefixRng.Find.Execute("Data di consegna: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($prefixFound) {
    $numStart = $prefixRng.End

    $numRng = $d.Range($numStart, $numStart + 2)

    if ($numRng.Text -eq "24") {
        # Toggling a direct-formatting property on just the "24" selection
        # before replacing its text forces the host to materialize the
        # replacement as its own run instead of silently folding it back
        # into the surrounding text - this is what gives us the three
        # separate runs seen in the target document.
        $numRng.Font.Bold = 1
        $numRng.Text = "31"

        $newNumRng = $d.Range($numStart, $numStart + 2)
        # Restore the normal (non-bold) formatting that matches the rest
        # of the line.
        $newNumRng.Font.Bold = 0
    }
}

# -------------------------------------------------------------------------
# Edit 2: remove "dashboard, " from the bullet point listing the app
# features, so the sentence reads "... (login, transazioni, analisi spese)."
# -------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Definire le funzionalità principali (login, dashboard, transazioni, analisi spese).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Definire le funzionalità principali (login, transazioni, analisi spese).",
    2
) | Out-Null
